# This script regenerates the save_data game log for jackson_luke (2021 season),
# replacing the old "Strike#" column with a recalculated "K" (strikeouts) column.
# The underlying simulation was re-run (new std/mean parameters, new s_vals draw),
# so the K values for each game row are recomputed/overwritten in-place.
#
# Column layout (row 1 headers): A=idx, B=date, C=TB, D=PC, E=dS0, F=dSF,
# G=K (strikeouts, the column being regenerated here), H=IP, I=I0, J=IF.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly simulated K values keyed by worksheet row number (rows 2-90 = games 0-88).
$kValues = @{
    2 = 1
    3 = 0
    4 = 3
    5 = 0
    6 = 0
    7 = 1
    8 = 1
    9 = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 3
    22 = 0
    23 = 3
    24 = 0
    25 = 4
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 3
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 3
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 2
    44 = 0
    45 = 0
    46 = 2
    47 = 1
    48 = 1
    49 = 1
    50 = 0
    51 = 1
    52 = 0
    53 = 2
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 1
    59 = 0
    60 = 0
    61 = 0
    62 = 1
    63 = 3
    64 = 1
    65 = 1
    66 = 1
    67 = 3
    68 = 0
    69 = 1
    70 = 2
    71 = 1
    72 = 0
    73 = 1
    74 = 0
    75 = 1
    76 = 0
    77 = 3
    78 = 0
    79 = 1
    80 = 0
    81 = 1
    82 = 0
    83 = 1
    84 = 1
    85 = 3
    86 = 1
    87 = 1
    88 = 0
    89 = 0
    90 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Output "Updated $($kValues.Count) K values in column G"
